$d = $word.ActiveDocument

$replacements = @(
    @{ Old = "2025-11-02 Sunday"; New = "2025-11-03 Monday" },
    @{ Old = "99÷6="; New = "63÷8=" },
    @{ Old = "12÷2="; New = "28÷4=" },
    @{ Old = "16÷6="; New = "57÷6=" },
    @{ Old = "99÷4="; New = "95÷5=" },
    @{ Old = "38÷2="; New = "36÷7=" },
    @{ Old = "72÷4="; New = "66÷2=" },
    @{ Old = "10÷7="; New = "76÷5=" },
    @{ Old = "51÷3="; New = "93÷5=" },
    @{ Old = "84÷5="; New = "84÷8=" },
    @{ Old = "56÷8="; New = "38÷8=" },
    @{ Old = "17÷9="; New = "34÷5=" },
    @{ Old = "87÷8="; New = "41÷9=" },
    @{ Old = "51÷8="; New = "76÷4=" },
    @{ Old = "25÷4="; New = "72÷7=" },
    @{ Old = "87÷5="; New = "50÷9=" },
    @{ Old = "54÷2="; New = "60÷3=" },
    @{ Old = "23÷8="; New = "15÷5=" },
    @{ Old = "33÷8="; New = "43÷5=" },
    @{ Old = "46÷7="; New = "30÷4=" },
    @{ Old = "44÷4="; New = "19÷2=" },
    @{ Old = "29÷6="; New = "12÷6=" },
    @{ Old = "83÷5="; New = "33÷5=" },
    @{ Old = "81÷3="; New = "18÷6=" },
    @{ Old = "57÷7="; New = "70÷8=" },
    @{ Old = "30÷7="; New = "26÷9=" }
)

foreach ($r in $replacements) {
    $d.Content.Find.Execute($r.Old, $true, $true, $false, $false, $false, $true, 1, $false, $r.New, 2)
}
